$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet: Total Trades 36 -> 37, Win Rate % 30.56 -> 29.73
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 37
$summary.Range("B9").Value = 29.73

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet: MarketMaking row - Trades 36 -> 37, Win Rate % 30.56 -> 29.73
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 37
$status.Range("G4").Value = 29.73

# ---------------------------------------------------------------------------
# 3) Append the newly-closed trade #37 as row 38 on both "All Trades" and
#    "MarketMaking" sheets (identical new trade record on each).
# ---------------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Range("A38").Value = 37

    # Date-looking text must stay literal text (not get auto-coerced into a
    # date serial number). Give it a text format, assign, then clear the
    # format again so no residual styling/number-format sticks to the cell.
    $ws.Range("B38").NumberFormat = "@"
    $ws.Range("B38").Value = "2026-02-17"
    $ws.Range("B38").ClearFormats()

    $ws.Range("C38").Value = "15:23:21"
    $ws.Range("D38").Value = "MarketMaking"
    $ws.Range("E38").Value = "UP"
    $ws.Range("F38").Value = 0.38
    $ws.Range("G38").Value = 0.38
    $ws.Range("H38").Value = "CLOSED"
    $ws.Range("I38").Value = 0
    $ws.Range("J38").Value = 0
    $ws.Range("K38").Value = 99.76000000000001
    $ws.Range("L38").Value = 0
    $ws.Range("M38").Value = 0
    $ws.Range("N38").Value = 0.6
    $ws.Range("O38").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P38").Value = "early_exit"
    $ws.Range("Q38").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
